# Revert "Updated report id and page code"
# This reverts the PBIReports sheet's PBIReportName/PBIReportPage columns
# back to their prior ("Financial overview" / ReportSection*) values, and
# restores the prior active-sheet / selection state.

$wb = $excel.ActiveWorkbook

$wsReports = $wb.Worksheets.Item("PBIReports")

# --- H column (PBIReportName): "Finance App" -> "Financial overview" for all data rows
$wsReports.Range("H2:H7").Value = "Financial overview"

# --- I column (PBIReportPage): restore the pre-update report-section identifiers
$wsReports.Range("I2").Value = "04fa320747962435bf38"
$wsReports.Range("I3").Value = "ReportSectionf72eb4d7e5e35db3b283"
$wsReports.Range("I4").Value = "ReportSection6a30609896651f006f0f"
$wsReports.Range("I5").Value = "ReportSection64d670dfa9da1a5b7033"
$wsReports.Range("I6").Value = "ReportSection6838cf9cda361d088e0a"
$wsReports.Range("I7").Value = "ReportSectionbb4917d9edb6d427282c"

# --- Restore prior view/selection state: PBIReports was the active tab with
# I1 selected; it reverts to being a background tab with C6 selected, while
# RCExtensionActions (selected D2:D7) becomes the active tab again.
$wsReports.Activate() | Out-Null
$wsReports.Range("C6").Select() | Out-Null

$wsActions = $wb.Worksheets.Item("RCExtensionActions")
$wsActions.Activate() | Out-Null
